$wb = $excel.ActiveWorkbook
$llm = $wb.ActiveSheet
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $llm)
$ws.Name = "VLM"

# Header row (copy style from LLM sheet header after setting values)
$ws.Range("A1").Value = 'Loadout'
$ws.Range("B1").Value = 'Scenario'
$ws.Range("C1").Value = 'Status'
$ws.Range("D1").Value = 'TTFT (s)'
$ws.Range("E1").Value = 'TPS'
$ws.Range("F1").Value = 'VRAM Peak (GB)'
$ws.Range("G1").Value = 'Text'
$llm.Range("A1:G1").Copy() | Out-Null
$ws.Range("A1:G1").PasteSpecial(-4122) | Out-Null

# Row 2
$ws.Range("A2").Value = 'base-qwen30-multi'
$ws.Range("B2").Value = 'bunny'
$ws.Range("C2").Value = 'PASSED'
$ws.Range("D2").Value = 0.3603103999994346
$ws.Range("E2").Value = 126.5041844137252
$ws.Range("F2").Value = 24.1025390625
$ws.Range("G2").Value = 'The animal in the video is a large, white, anthropomorphic rabbit. It is shown in a series of scenes in a lush, green field.
- In the first scene, the rabbit is lying on its back on the grass, with a large, pink butterfly resting on its belly.
- In the second scene, the rabbit is sitting up, looking down at the butterfly, which is now perched on its head.
- In the final scene, the rabbit is standing in a field, looking down at a small, red, round object on the ground, possibly a flower or a ball.'

# Row 3
$ws.Range("A3").Value = 'base-qwen30-multi'
$ws.Range("B3").Value = 'jarvis_logo'
$ws.Range("C3").Value = 'PASSED'
$ws.Range("D3").Value = 0.1563606999989133
$ws.Range("E3").Value = 176.0141120761007
$ws.Range("F3").Value = 24.1025390625
$ws.Range("G3").Value = 'This image features a stylized, futuristic robot head, which appears to be a representation of Iron Man''s helmet. The helmet is centrally positioned within a circular, high-tech interface that resembles a heads-up display (HUD) or targeting system.
Key visual elements include:
- A glowing blue, metallic helmet with a sleek, angular design.
- Bright, glowing green eyes that give it an intense, focused look.
- Red accents on the sides of the helmet, possibly representing earpieces or sensors.
- A circular frame surrounding the helmet, composed of concentric rings with glowing blue lines and markings, suggesting a digital or scanning interface.
- The overall color scheme is dominated by shades of blue and black, creating a dark, high-tech, and cybernetic atmosphere.
The image is likely a logo or icon for a game, app, or digital platform related to Iron Man or a similar sci-fi theme.'

# Row 4
$ws.Range("A4").Value = 'base-qwen30-multi'
$ws.Range("B4").Value = 'three_objects'
$ws.Range("C4").Value = 'PASSED'
$ws.Range("D4").Value = 0.1310329000043566
$ws.Range("E4").Value = 146.9501913408375
$ws.Range("F4").Value = 24.1025390625
$ws.Range("G4").Value = 'Based on the image provided, we can identify the following:
- There is one red circle in the upper left.
- There is one blue circle in the upper right.
- There is one green circle in the lower center.
Each of these is a distinct, colored circle.
Therefore, there are 3 colored circles in the photo.'

# Row 5
$ws.Range("A5").Value = 'base-qwen30-multi'
$ws.Range("B5").Value = 'traffic'
$ws.Range("C5").Value = 'PASSED'
$ws.Range("D5").Value = 1.030388999999559
$ws.Range("E5").Value = 39.32765441986089
$ws.Range("F5").Value = 24.1025390625
$ws.Range("G5").Value = 'Based on the video frames provided, the following types of vehicles are visible:
- **Car**: A white car is seen driving into the parking lot and then parking.
- **Bicycle**: A person is seen riding a bicycle across the parking lot.'

